$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 6.305
$ws.Range("A12").Value = -21.489
$ws.Range("B23").Value = 8.472
$ws.Range("C24").Value = -12.572
$ws.Range("B28").Value = 5.142
$ws.Range("A32").Value = -21.243
$ws.Range("B32").Value = 7.063
$ws.Range("B34").Value = 6.952000000000001
$ws.Range("A36").Value = -20.724
$ws.Range("A38").Value = -20.354
$ws.Range("C38").Value = -11.836
$ws.Range("B42").Value = 9.103999999999999
$ws.Range("A46").Value = -21.728
$ws.Range("C52").Value = -11.701
$ws.Range("A54").Value = -21.37
$ws.Range("B54").Value = 5.406000000000001
$ws.Range("A55").Value = -22.184
$ws.Range("A67").Value = -21.531
$ws.Range("A69").Value = -21.422
$ws.Range("A72").Value = -21.695
$ws.Range("C78").Value = -12.539
$ws.Range("C83").Value = -13.409
$ws.Range("C85").Value = -12.335
$ws.Range("C86").Value = -13.725
$ws.Range("A91").Value = -20.811
$ws.Range("C96").Value = -10.156
$ws.Range("B97").Value = 5.391000000000001
$ws.Range("A99").Value = -21.282
$ws.Range("B99").Value = 5.671000000000001
$ws.Range("B101").Value = 5.496
$ws.Range("C103").Value = -12.411
$ws.Range("A104").Value = -21.437
